$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "44.066.86"
$ws.Range("D3").Value = "2.359.10"
$ws.Range("E3").Value = "  +0.48%  "
$ws.Range("E4").Value = "  +0.15%  "
$ws.Range("D5").Value = "'0.695"
$ws.Range("E5").Value = "  +5.65%  "
$ws.Range("D6").Value = "'241.79"
$ws.Range("E6").Value = "  +2.85%  "
$ws.Range("D7").Value = "'76.59"
$ws.Range("E7").Value = "  +4.50%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("D9").Value = "'0.627"
$ws.Range("E9").Value = "  +20.05%  "
$ws.Range("E10").Value = "  +3.88%  "
$ws.Range("E11").Value = "  +0.82%  "
$ws.Range("D12").Value = "'33.32"
$ws.Range("E12").Value = "  +21.18%  "
$ws.Range("D13").Value = "'7.49"
$ws.Range("E13").Value = "  +12.29%  "
$ws.Range("D14").Value = "'0.109"
$ws.Range("E14").Value = "  +1.96%  "
$ws.Range("D15").Value = "2.708.49"
$ws.Range("E15").Value = "  +0.43%  "
$ws.Range("D16").Value = "'16.75"
$ws.Range("E16").Value = "  -2.80%  "
$ws.Range("D17").Value = "'0.928"
$ws.Range("E17").Value = "  +5.69%  "
$ws.Range("D18").Value = "2.355.27"
$ws.Range("E18").Value = "  -0.30%  "
$ws.Range("D19").Value = "43.968.71"
$ws.Range("E19").Value = "  +0.79%  "
$ws.Range("E20").Value = "  +2.64%  "
$ws.Range("D21").Value = "'6.70"
$ws.Range("E21").Value = "  +4.93%  "
$ws.Range("D22").Value = "'77.93"
$ws.Range("E22").Value = "  +2.70%  "
$ws.Range("D23").Value = "'261.49"
$ws.Range("E23").Value = "  +4.29%  "
$ws.Range("E24").Value = "  -0.09%  "
$ws.Range("D25").Value = "'2.54"
$ws.Range("E25").Value = "  +2.48%  "
$ws.Range("D26").Value = "'3.64"
$ws.Range("E26").Value = "  -4.51%  "
$ws.Range("E27").Value = "  +18.12%  "
$ws.Range("D28").Value = "'10.95"
$ws.Range("E28").Value = "  +6.73%  "
$ws.Range("E29").Value = "  +4.04%  "
$ws.Range("D30").Value = "'23.23"
$ws.Range("E30").Value = "  +3.49%  "
$ws.Range("D31").Value = "'175.13"
$ws.Range("E31").Value = "  +1.63%  "
$ws.Range("E32").Value = "  -4.08%  "
$ws.Range("D33").Value = "'0.137"
$ws.Range("E33").Value = "  +5.40%  "
$ws.Range("D34").Value = "'5.41"
$ws.Range("E34").Value = "  +6.58%  "
$ws.Range("E35").Value = "  +9.11%  "
$ws.Range("E36").Value = "  +7.38%  "
$ws.Range("D37").Value = "'3.83"
$ws.Range("E37").Value = "  +2.02%  "
$ws.Range("E38").Value = "  +0.15%  "
$ws.Range("D39").Value = "'6.42"
$ws.Range("E39").Value = "  -0.17%  "
$ws.Range("D40").Value = "'0.0284"
$ws.Range("E40").Value = "  +7.29%  "
$ws.Range("D41").Value = "'0.219"
$ws.Range("E41").Value = "  +21.18%  "
$ws.Range("B42").Value = "Cronos"
$ws.Range("C42").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D42").Value = "'0.108"
$ws.Range("E42").Value = "  +12.53%  "
$ws.Range("B43").Value = "InjectiveProtocol"
$ws.Range("C43").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D43").Value = "'19.31"
$ws.Range("E43").Value = "  -0.79%  "
$ws.Range("D44").Value = "'9.17"
$ws.Range("E44").Value = "  +3.25%  "
$ws.Range("E45").Value = "  -0.13%  "
$ws.Range("E46").Value = "  +6.16%  "
$ws.Range("E47").Value = "  +10.06%  "
$ws.Range("E48").Value = "  +3.97%  "
$ws.Range("D49").Value = "'102.27"
$ws.Range("E49").Value = "  +3.53%  "
$ws.Range("E50").Value = "  +2.81%  "
$ws.Range("D51").Value = "'56.56"
$ws.Range("E51").Value = "  +10.80%  "
